$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.048.96"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "1.647.63"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'215.37"
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("D6").Value = "'0.5092"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.2581"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").Value = "'0.06416"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").Value = "'19.58"
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.651.19"
$ws.Range("E12").Value = "  -3.43%  "
$ws.Range("D13").Value = "'4.260"
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("D14").Value = "1.875.75"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").Value = "'0.5454"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").Value = "0.0₅7972"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "'63.73"
$ws.Range("E17").Value = "  -4.58%  "
$ws.Range("D18").Value = "26.069.34"
$ws.Range("E18").Value = "  -4.18%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'206.59"
$ws.Range("E20").Value = "  -4.89%  "
$ws.Range("D21").Value = "'4.372"
$ws.Range("E21").Value = "  -4.56%  "
$ws.Range("D22").Value = "'10.01"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "'5.996"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'1.867"
$ws.Range("E25").Value = "  +7.30%  "
$ws.Range("D26").Value = "'143.04"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "'0.1165"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "'6.893"
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D29").Value = "'15.77"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "'0.05066"
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("D31").Value = "'1.240"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "'3.316"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").Value = "'3.221"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "'1.543"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("D35").Value = "'2.346"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("D36").Value = "'0.9126"
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("D37").Value = "'2.646"
$ws.Range("E37").Value = "  -6.07%  "
$ws.Range("D38").Value = "'0.5687"
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("D39").Value = "1.147.44"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("D40").Value = "'0.01569"
$ws.Range("E40").Value = "  -3.91%  "
$ws.Range("D41").Value = "'2.565"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'5.653"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "'0.8228"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "'99.69"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "1.787.86"
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("D48").Value = "'0.4532"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "'55.07"
$ws.Range("E50").Value = "  -3.35%  "
$ws.Range("D51").Value = "'7.828"
$ws.Range("E51").Value = "  -2.83%  "
